# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" columns (F/G) for the
# zh-cn and de-de handback rows, refreshes the handback timestamp/status strings.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: status column values (zh-cn / de-de) move from "Ready for
# handoff" to "Handed back: in sync with en-US" for both files.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusText
$wsOverview.Range("C2").Value = $statusText
$wsOverview.Range("B3").Value = $statusText
$wsOverview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------------
# Helper that rewrites a language sheet (zh-cn / de-de) with the new
# "Latest Target File" (F) / "Latest Handback File" (G) hyperlinks, and
# refreshes the status + handback datetime columns.
# ---------------------------------------------------------------------------
function Update-HandbackSheet {
    param($ws, [string]$lang, [string]$xlfName, [string]$handoffSha, [string]$handbackSha, [string]$handbackDateTime)

    $mdName = "807c49aa-bb0d-4129-a629-50310e42ed4f.md"
    $ffName = "ffff08f1150c-0934-43ad-99f7-f14e0b137dbd.md"

    $ffUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5d5e83cdaf46ab7aa1da116fc6ffbd75324302a9/e2e/$ffName"
    $targetUrl = "https://github.com/OpenLocalizationTestOrg/oltest.$lang/blob/$handbackSha/e2e/$mdName"
    $handbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$handbackSha/ol-handback/OpenLocalizationTestOrg/oltest.$lang/ci/hb/$xlfName"

    # Status column refreshed to reflect the handback completion.
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Handback DateTime column.
    $ws.Range("H2").Value = $handbackDateTime
    $ws.Range("H3").Value = $handbackDateTime

    # Remove & recreate the row-3 hyperlinks so new relationship ids slot in
    # after the newly-added row-2 hyperlinks (mirrors how the report
    # generator lays the <hyperlinks> collection out).
    function Remove-HL($addr) {
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Address() -eq $addr) {
                $hl.Delete()
                return
            }
        }
    }

    Remove-HL '$A$3'
    Remove-HL '$B$3'
    Remove-HL '$D$3'

    # New row 2: Latest Target File (F2) / Latest Handback File (G2)
    $ws.Hyperlinks.Add($ws.Range("F2"), $targetUrl, "", "", $mdName)
    $ws.Hyperlinks.Add($ws.Range("G2"), $handbackUrl, "", "", $xlfName)

    # Re-add row 3 original hyperlinks (A3/B3/D3) unchanged in content.
    $ws.Hyperlinks.Add($ws.Range("A3"), $ffUrl, "", "", $ffName)
    $ws.Hyperlinks.Add($ws.Range("B3"), $ffUrl, "", "", ".md")
    $ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$handoffSha/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/ci/ht/$xlfName", "", "", $xlfName)

    # New row 3: Latest Target File (F3) / Latest Handback File (G3)
    $ws.Hyperlinks.Add($ws.Range("F3"), $targetUrl, "", "", $mdName)
    $ws.Hyperlinks.Add($ws.Range("G3"), $handbackUrl, "", "", $xlfName)

    # All hyperlink cells should carry the same "HyperLink" cell style.
    $ws.Range("A3").Style = "HyperLink"
    $ws.Range("B3").Style = "HyperLink"
    $ws.Range("D3").Style = "HyperLink"
    $ws.Range("F2").Style = "HyperLink"
    $ws.Range("G2").Style = "HyperLink"
    $ws.Range("F3").Style = "HyperLink"
    $ws.Range("G3").Style = "HyperLink"
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-HandbackSheet $wsZhCn "zh-cn" `
    "807c49aa-bb0d-4129-a629-50310e42ed4f.c9f3bf6d024ae69c20d221e7317f1a9ee5de055e.zh-cn.xlf" `
    "4ecf118ea642c17f1a3079ea07b65a5777e6f494" `
    "c9f3bf6d024ae69c20d221e7317f1a9ee5de055e" `
    "2016-03-12 00:44:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
Update-HandbackSheet $wsDeDe "de-de" `
    "807c49aa-bb0d-4129-a629-50310e42ed4f.c9f3bf6d024ae69c20d221e7317f1a9ee5de055e.de-de.xlf" `
    "926b9bce2fba568c0257d2adf822a88960f23c22" `
    "c9f3bf6d024ae69c20d221e7317f1a9ee5de055e" `
    "2016-03-12 00:44:33"
